$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.693.63"
$ws.Range("E2").Value = "  +4.22%  "
$ws.Range("D3").Value = "2.258.07"
$ws.Range("E3").Value = "  +2.48%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.03"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "91.21"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.68%  "
$ws.Range("E7").Value = "  +3.46%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.480"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.10"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.19"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.43%  "
$ws.Range("E12").Value = "  +2.03%  "
$ws.Range("E13").Value = "  +1.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.57"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.99%  "
$ws.Range("D15").Value = "2.607.74"
$ws.Range("E15").Value = "  +2.44%  "
$ws.Range("E16").Value = "  +2.61%  "
$ws.Range("D17").Value = "2.260.53"
$ws.Range("E17").Value = "  -4.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.758"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.27%  "
$ws.Range("D19").Value = "41.601.09"
$ws.Range("E19").Value = "  +4.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.44"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +11.01%  "
$ws.Range("E21").Value = "  +1.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.90"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.54"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "240.53"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.42%  "
$ws.Range("E25").Value = "  +4.40%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  +5.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.02"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.48"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.11%  "
$ws.Range("E30").Value = "  -0.96%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "160.90"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "34.44"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.14"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0742"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.23%  "
$ws.Range("E36").Value = "  -0.41%  "
$ws.Range("E37").Value = "  +2.31%  "
$ws.Range("E38").Value = "  +2.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.60"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.01%  "
$ws.Range("E40").Value = "  +4.18%  "
$ws.Range("E41").Value = "  +3.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.91"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.05%  "
$ws.Range("D43").Value = "2.058.59"
$ws.Range("E43").Value = "  -0.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.69"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.13%  "
$ws.Range("E45").Value = "  +2.88%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.14"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.81%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.86"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.86%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.04"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.89%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.51"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.78%  "
$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "72.33"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.22%  "
$ws.Range("E51").Value = "  +2.70%  "
